$d = $word.ActiveDocument

# Locate the "Requisitos" bullet-list paragraph: it's the paragraph that
# contains the "LOM3206" requirement line (together with LOM3215, LOM3231,
# LOM3234 as separate runs, each followed by a manual line break).
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*LOM3206*LOM3231*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the Requisitos paragraph"
}

# Rebuild that paragraph's runs with "LOM3231 - ... (Indicação de Conjunto)"
# moved so it comes first (before the LOM3206 line), keeping every
# requirement line as its own run ending in a manual line break, exactly
# like the rest of the paragraph's runs already do.
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
    '<w:r><w:t>LOM3231 -  M' + [char]0xE9 + 'todos Experimentais da F' + [char]0xED + 'sica IV  (Indica' + [char]0xE7 + [char]0xE3 + 'o de Conjunto)</w:t><w:br/></w:r>' +
    '<w:r><w:t>LOM3206 -  Eletr' + [char]0xF4 + 'nica  (Requisito)</w:t><w:br/></w:r>' +
    '<w:r><w:t>LOM3215 -  F' + [char]0xED + 'sica do Estado S' + [char]0xF3 + 'lido  (Requisito)</w:t><w:br/></w:r>' +
    '<w:r><w:t>LOM3234 -  ' + [char]0xD3 + 'ptica F' + [char]0xED + 'sica  (Requisito)</w:t><w:br/></w:r>' +
    '</w:p>'

[void]$target.Range.InsertXML($newXml)
